$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1): add a handoff-report row for the new file
# 413ded0f-c5bb-4c26-943b-fce76493f7f9.md
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Insert a new row 3 - this copies down the formatting (styles) from row 2,
# matching the existing "Ready for handoff" rows.
$wsOverview.Rows(3).Insert()

$wsOverview.Range("A3").Value = "413ded0f-c5bb-4c26-943b-fce76493f7f9.md"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-23 08:41:25"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0a258592573329ff777d89e0173949bf1577b6c8/e2e/413ded0f-c5bb-4c26-943b-fce76493f7f9.md", "", "", "413ded0f-c5bb-4c26-943b-fce76493f7f9.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2): add the matching handoff row for the new file
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Rows(3).Insert()

$wsZhCn.Range("A3").Value = "413ded0f-c5bb-4c26-943b-fce76493f7f9.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "413ded0f-c5bb-4c26-943b-fce76493f7f9.6f6ef1c0f04faa52a3b826b4db3169c71ecf7fcb.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-23 08:41:21"
$wsZhCn.Range("H3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("J3").Value = "Include"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0a258592573329ff777d89e0173949bf1577b6c8/e2e/413ded0f-c5bb-4c26-943b-fce76493f7f9.md", "", "", "413ded0f-c5bb-4c26-943b-fce76493f7f9.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/561769c95e1f15ca3e2a640f6d5766b0468f0b73/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/413ded0f-c5bb-4c26-943b-fce76493f7f9.6f6ef1c0f04faa52a3b826b4db3169c71ecf7fcb.zh-cn.xlf", "", "", "413ded0f-c5bb-4c26-943b-fce76493f7f9.6f6ef1c0f04faa52a3b826b4db3169c71ecf7fcb.zh-cn.xlf")

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3): add the matching handoff row for the new file
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Rows(3).Insert()

$wsDeDe.Range("A3").Value = "413ded0f-c5bb-4c26-943b-fce76493f7f9.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "413ded0f-c5bb-4c26-943b-fce76493f7f9.6f6ef1c0f04faa52a3b826b4db3169c71ecf7fcb.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-23 08:41:25"
$wsDeDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("J3").Value = "Include"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0a258592573329ff777d89e0173949bf1577b6c8/e2e/413ded0f-c5bb-4c26-943b-fce76493f7f9.md", "", "", "413ded0f-c5bb-4c26-943b-fce76493f7f9.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ac8877ccb63151644a7a6d4dd2c4feeedaee0cfb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/413ded0f-c5bb-4c26-943b-fce76493f7f9.6f6ef1c0f04faa52a3b826b4db3169c71ecf7fcb.de-de.xlf", "", "", "413ded0f-c5bb-4c26-943b-fce76493f7f9.6f6ef1c0f04faa52a3b826b4db3169c71ecf7fcb.de-de.xlf")
